$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 17.306265
$ws.Cells.Item(2,8).Value = 51.918795
$ws.Cells.Item(2,9).Value = 0.5463168539988408
$ws.Cells.Item(2,10).Value = 0.5463168539988407
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 14.792724
$ws.Cells.Item(2,14).Value = 44.378172
$ws.Cells.Item(2,15).Value = 0.0761423615956231
$ws.Cells.Item(2,16).Value = 0.0761423615956231
$ws.Cells.Item(2,17).Value = 256.00680161586
$ws.Cells.Item(2,18).Value = 2304.06121454274
$ws.Cells.Item(2,19).Value = 0.04159785544296297
$ws.Cells.Item(2,20).Value = 0.04159785544296295

$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 17.306265
$ws.Cells.Item(3,8).Value = 51.918795
$ws.Cells.Item(3,9).Value = 0.5463168539988408
$ws.Cells.Item(3,10).Value = 0.5463168539988407
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 44.28072733333334
$ws.Cells.Item(3,14).Value = 132.842182
$ws.Cells.Item(3,15).Value = 0.2279255093471533
$ws.Cells.Item(3,16).Value = 0.2279255093471533
$ws.Cells.Item(3,17).Value = 766.3340016234101
$ws.Cells.Item(3,18).Value = 6897.006014610691
$ws.Cells.Item(3,19).Value = 0.1245195472126202
$ws.Cells.Item(3,20).Value = 0.1245195472126201

$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 17.306265
$ws.Cells.Item(4,8).Value = 51.918795
$ws.Cells.Item(4,9).Value = 0.5463168539988408
$ws.Cells.Item(4,10).Value = 0.5463168539988407
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 135.203738
$ws.Cells.Item(4,14).Value = 405.611214
$ws.Cells.Item(4,15).Value = 0.6959321290572236
$ws.Cells.Item(4,16).Value = 0.6959321290572236
$ws.Cells.Item(4,17).Value = 2339.87171881857
$ws.Cells.Item(4,18).Value = 21058.84546936713
$ws.Cells.Item(4,19).Value = 0.3801994513432576
$ws.Cells.Item(4,20).Value = 0.3801994513432576

$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 3.344413333333334
$ws.Cells.Item(5,8).Value = 10.03324
$ws.Cells.Item(5,9).Value = 0.1055750256186672
$ws.Cells.Item(5,10).Value = 0.1055750256186672
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 14.792724
$ws.Cells.Item(5,14).Value = 44.378172
$ws.Cells.Item(5,15).Value = 0.0761423615956231
$ws.Cells.Item(5,16).Value = 0.0761423615956231
$ws.Cells.Item(5,17).Value = 49.47298338192
$ws.Cells.Item(5,18).Value = 445.25685043728
$ws.Cells.Item(5,19).Value = 0.008038731776123728
$ws.Cells.Item(5,20).Value = 0.008038731776123726

$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 3.344413333333334
$ws.Cells.Item(6,8).Value = 10.03324
$ws.Cells.Item(6,9).Value = 0.1055750256186672
$ws.Cells.Item(6,10).Value = 0.1055750256186672
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 44.28072733333334
$ws.Cells.Item(6,14).Value = 132.842182
$ws.Cells.Item(6,15).Value = 0.2279255093471533
$ws.Cells.Item(6,16).Value = 0.2279255093471533
$ws.Cells.Item(6,17).Value = 148.0930549032978
$ws.Cells.Item(6,18).Value = 1332.83749412968
$ws.Cells.Item(6,19).Value = 0.02406324148847347
$ws.Cells.Item(6,20).Value = 0.02406324148847347

$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 3.344413333333334
$ws.Cells.Item(7,8).Value = 10.03324
$ws.Cells.Item(7,9).Value = 0.1055750256186672
$ws.Cells.Item(7,10).Value = 0.1055750256186672
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 135.203738
$ws.Cells.Item(7,14).Value = 405.611214
$ws.Cells.Item(7,15).Value = 0.6959321290572236
$ws.Cells.Item(7,16).Value = 0.6959321290572236
$ws.Cells.Item(7,17).Value = 452.1771840837068
$ws.Cells.Item(7,18).Value = 4069.594656753361
$ws.Cells.Item(7,19).Value = 0.07347305235406996
$ws.Cells.Item(7,20).Value = 0.07347305235406996

$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 11.02739433333333
$ws.Cells.Item(8,8).Value = 33.082183
$ws.Cells.Item(8,9).Value = 0.3481081203824922
$ws.Cells.Item(8,10).Value = 0.3481081203824921
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 14.792724
$ws.Cells.Item(8,14).Value = 44.378172
$ws.Cells.Item(8,15).Value = 0.0761423615956231
$ws.Cells.Item(8,16).Value = 0.0761423615956231
$ws.Cells.Item(8,17).Value = 163.125200812164
$ws.Cells.Item(8,18).Value = 1468.126807309476
$ws.Cells.Item(8,19).Value = 0.02650577437653641
$ws.Cells.Item(8,20).Value = 0.02650577437653641

$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 11.02739433333333
$ws.Cells.Item(9,8).Value = 33.082183
$ws.Cells.Item(9,9).Value = 0.3481081203824922
$ws.Cells.Item(9,10).Value = 0.3481081203824921
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 44.28072733333334
$ws.Cells.Item(9,14).Value = 132.842182
$ws.Cells.Item(9,15).Value = 0.2279255093471533
$ws.Cells.Item(9,16).Value = 0.2279255093471533
$ws.Cells.Item(9,17).Value = 488.3010416714785
$ws.Cells.Item(9,18).Value = 4394.709375043306
$ws.Cells.Item(9,19).Value = 0.07934272064605968
$ws.Cells.Item(9,20).Value = 0.07934272064605967

$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 11.02739433333333
$ws.Cells.Item(10,8).Value = 33.082183
$ws.Cells.Item(10,9).Value = 0.3481081203824922
$ws.Cells.Item(10,10).Value = 0.3481081203824921
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 135.203738
$ws.Cells.Item(10,14).Value = 405.611214
$ws.Cells.Item(10,15).Value = 0.6959321290572236
$ws.Cells.Item(10,16).Value = 0.6959321290572236
$ws.Cells.Item(10,17).Value = 1490.944934266685
$ws.Cells.Item(10,18).Value = 13418.50440840016
$ws.Cells.Item(10,19).Value = 0.2422596253598961
$ws.Cells.Item(10,20).Value = 0.242259625359896

